$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Metadata" (sheet1) ---

# Remove the duplicate "Contact" row (old row 11 was an exact duplicate of row 10)
$ws1.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> Alvearie Team
$ws1.Range("B9").Value = "Alvearie Team"

# The old "Contact" / "No display for ContactDetail" row (now row 10 after the
# delete above) becomes "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# --- Sheet "Elements" (sheet2) ---
# Short/Definition columns for the root Extension row now read "Measure Weight"
$ws2.Range("K2").Value = "Measure Weight"
$ws2.Range("L2").Value = "Measure Weight"
